$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select rows 2-4 (RefID 1, 4, 10 - submissions without a PMID) and delete them entirely,
# shifting the remaining rows up.
$rng = $ws.Range("A2:XFD4")
$rng.Select()
$rng.EntireRow.Delete()
